$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking strings are not
# auto-converted to numbers by Excel, matching the inlineStr cells in the source.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "28.364.34"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.574.46"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "212.06"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "0.489"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "44.58"
$ws.Range("E8").Value = "  -4.43%  "
$ws.Range("D9").Value = "23.79"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").Value = "0.0896"
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").Value = "1.798.23"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").Value = "1.567.70"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("E16").Value = "  -1.39%  "
$ws.Range("D17").Value = "28.366.64"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "61.58"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").Value = "230.32"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").Value = "7.43"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "3.96"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").Value = "2.03"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "151.42"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").Value = "14.93"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "0.104"
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "6.36"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  +3.13%  "
$ws.Range("E32").Value = "  -3.78%  "
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("D34").Value = "3.10"
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("D35").Value = "1.386.61"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("E36").Value = "  +5.45%  "
$ws.Range("E37").Value = "  -3.06%  "
$ws.Range("D38").Value = "2.36"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").Value = "0.517"
$ws.Range("E41").Value = "  -2.63%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "1.90"
$ws.Range("E43").Value = "  +2.40%  "
$ws.Range("D44").Value = "0.786"
$ws.Range("E44").Value = "  -1.24%  "
$ws.Range("D45").Value = "0.0467"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("E46").Value = "  -4.10%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "62.46"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "0.920"
$ws.Range("E48").Value = "  -6.02%  "
$ws.Range("D49").Value = "1.710.98"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").Value = "85.46"
$ws.Range("E51").Value = "  -0.49%  "

# Restore default (Normal) style on column D so no stray style index is left on cells
$priceCol.Style = "Normal"

